$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) LiaJourney sheet: small wording tweaks
# ---------------------------------------------------------------------------
$wsJourney = $wb.Worksheets.Item("LiaJourney")
$wsJourney.Range("B5").Value = "wait"
$wsJourney.Range("B11").Value = "I can't resist you anymore"

# ---------------------------------------------------------------------------
# 2) Split "cumcontrol" into "cumcontrol1" + "cumcontrol2", and restore a
#    standalone "dickpic" sheet with its original content (QA round 2:
#    diversification of the cum-control script + a dedicated dick pic sheet).
# ---------------------------------------------------------------------------
$wsCumControl = $wb.Worksheets.Item("cumcontrol")
$wsDickPic = $wb.Worksheets.Item("dickpic")

# Duplicate the current "dickpic" sheet and place the copy right after
# "cumcontrol". That duplicate will become the new "cumcontrol2" sheet, while
# the original "dickpic" sheet (further down) keeps its original content and
# simply gets shifted after "cumcontrol2".
$cumControlIndex = $wsCumControl.Index
$wsDickPic.Copy($null, $wsCumControl)
$wsCumControl2 = $wb.Worksheets.Item($cumControlIndex + 1)
$wsCumControl2.Name = "cumcontrol2"

# Rename the original "cumcontrol" sheet to "cumcontrol1"
$wsCumControl.Name = "cumcontrol1"

# ---- cumcontrol1: reword the existing delay/sync/edge variant-1 lines ----
$wsCumControl.Range("B2").Value = "if you finish before you see what I'm sending next you'll regret it"

$wsCumControl.Range("B3").Value = "wait wait wait... I have one more thing for you before you finish"
$wsCumControl.Range("C3").Value = "DELAY. Send final PPV."

$wsCumControl.Range("B4").Value = "I want to feel it at the same time love... watch this first"
$wsCumControl.Range("C4").Value = "SYNC variant. Send PPV."

$wsCumControl.Range("B5").Value = "okay NOW we can go together... open this"
$wsCumControl.Range("C5").Value = "SYNC. Send PPV."

$wsCumControl.Range("B6").Value = "you better not be close already... I have more to show you"

$wsCumControl.Range("B7").Value = "not yet... I said not yet love"
$wsCumControl.Range("C7").Value = "CONTROL. More PPVs to send. Create urgency to open next."

# ---- cumcontrol2: replace the duplicated "dickpic" rows with a brand new
#      delay/sync/edge variant-2 script ----
$wsCumControl2.Range("A2").Value = "delay2"
$wsCumControl2.Range("B2").Value = "hold on just a little longer, I promise this next one is worth it"
$wsCumControl2.Range("C2").Value = "DELAY variant."

$wsCumControl2.Range("A3").Value = "delay1"
$wsCumControl2.Range("B3").Value = "don't you dare... not until you see what I just did"
$wsCumControl2.Range("C3").Value = "DELAY. Send PPV."

$wsCumControl2.Range("A4").Value = "sync2"
$wsCumControl2.Range("B4").Value = "let's do this together love... but you have to open this first"
$wsCumControl2.Range("C4").Value = "SYNC variant."

$wsCumControl2.Range("A5").Value = "sync1"
$wsCumControl2.Range("B5").Value = "okay I'm ready now too... watch this with me"
$wsCumControl2.Range("C5").Value = "SYNC. Send PPV."

$wsCumControl2.Range("A6").Value = "edge2"
$wsCumControl2.Range("B6").Value = "patience... the best part hasn't even happened yet"
$wsCumControl2.Range("C6").Value = "EDGE variant."

$wsCumControl2.Range("A7").Value = "edge1"
$wsCumControl2.Range("B7").Value = "slow down love... I'm not letting you off that easy"
$wsCumControl2.Range("C7").Value = "CONTROL."
